# Updated cryptos list on Mon Aug 12 15:42:56 UTC 2024 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns with latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.476.06'
$ws.Range("E2").Value = '  -1.23%  '

$ws.Range("D3").Value = '2.638.56'
$ws.Range("E3").Value = '  +1.01%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.47%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.994'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.49%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.572'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.44%  '

$ws.Range("D9").Value = '2.655.88'
$ws.Range("E9").Value = '  +1.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.48'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.107'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.340'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("E13").Value = '  -1.42%  '

$ws.Range("D14").Value = '3.099.83'
$ws.Range("E14").Value = '  +1.04%  '

$ws.Range("D15").Value = '59.357.89'
$ws.Range("E15").Value = '  -1.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000139'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.44%  '

$ws.Range("D18").Value = '2.692.01'
$ws.Range("E18").Value = '  +3.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.62'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '345.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.26%  '

$ws.Range("E23").Value = '  +0.85%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.78'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.425'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.42%  '

$ws.Range("D26").Value = '2.765.90'
$ws.Range("E26").Value = '  +1.44%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.991'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.161'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.99%  '

$ws.Range("D29").Value = '0.0₃0824'
$ws.Range("E29").Value = '  +1.98%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.37%  '

$ws.Range("E31").Value = '  -0.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.54'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.54%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.47%  '

$ws.Range("E34").Value = '  -0.90%  '

$ws.Range("E35").Value = '  +16.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '149.44'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.00%  '

$ws.Range("E38").Value = '  +1.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.872'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.65'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.87%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.75'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '284.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.616'
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0995'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.992'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.80'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0545'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0233'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.73%  '

$ws.Range("E51").Value = '  -1.37%  '
